# New column C: "GPU V2" timings, added alongside the existing CPU VERSION / GPU V1 columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "GPU V2"

$ws.Range("C2").Value = 24.364999999999998
$ws.Range("C3").Value = 20.271000000000001
$ws.Range("C4").Value = 15.59
$ws.Range("C5").Value = 19.966999000000001
$ws.Range("C6").Value = 15.97

# C5 keeps the same custom numeric format already used by column B.
$ws.Range("C5").NumberFormat = "#,##0.000000"

# Widen column C to fit its new contents.
$ws.Columns("C:C").AutoFit()

# Commit-message note, placed under the new data (C7).
$ws.Range("C7").Value = "add constant memory for rowPointers and destinations"

# New column D: header for the next version (GPU v3).
$ws.Range("D1").Value = "GPU v3"

$ws.Range("H13").Select()
